# Applies the "Presentacion comentarios fb page" edit to Sheet1.
# - Swaps the E1/F1 header labels ("Comments" <-> "Id").
# - Replaces the content of data rows 2 and 3 (new tweet text / label / rate /
#   numeric id in E / cleared F / updated Neu_com).
# - Removes the old data rows 4-7 so only the header + 2 data rows remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap E1 / F1 labels -------------------------------------
$ws.Range("E1").Value = "Id"
$ws.Range("F1").Value = "Comments"

# --- Row 2 ----------------------------------------------------------------
$ws.Range("B2").Value = "GRACIAS POR ESTAR EN EL DIRECTO CON NEYMAR, NINJA Y FLAKES`nVOY PARA VER EL ÚLTIMO COMBATE`nOS QUIERO ❤️ https://t.co/IVmdUiMTdU"
$ws.Range("C2").Value = "Neutral"
$ws.Range("D2").Value = 0.1213210496251103
$ws.Range("E2").Value = 1397656014909026000
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

# --- Row 3 ----------------------------------------------------------------
$ws.Range("B3").Value = "Stream con Neymar en 5 mins.`n❤️💛 https://t.co/hWQuFo0UxU https://t.co/6zdLLHwpCB"
$ws.Range("C3").Value = "Positive"
$ws.Range("D3").Value = 0.9041762621135991
$ws.Range("E3").Value = 1397624785585193000
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Re-fit rows 2/3 so the multi-line text doesn't leave a stray custom row
# height behind (matches the source file, which has no explicit row heights).
$ws.Rows("2:3").AutoFit()

# --- Remove old rows 4-7 so the sheet ends at row 3 ------------------------
$ws.Rows("4:7").Delete()
